{"js": "// Replace the three-digit-by-one-digit multiplication problems with new\n// operands, keeping every other part of the document (layout, fonts,\n// formatting, trailing \"=\") unchanged. Each original expression is unique\n// in the document, so a plain search + full-text replace per pair is safe.\nconst replacements = [\n  [\"789\u00d74=\", \"217\u00d76=\"],\n  [\"362\u00d73=\", \"619\u00d75=\"],\n  [\"791\u00d74=\", \"111\u00d78=\"],\n  [\"201\u00d79=\", \"404\u00d72=\"],\n  [\"329\u00d78=\", \"363\u00d76=\"],\n  [\"262\u00d77=\", \"547\u00d79=\"],\n  [\"225\u00d77=\", \"310\u00d78=\"],\n  [\"443\u00d73=\", \"315\u00d75=\"],\n  [\"867\u00d76=\", \"955\u00d77=\"],\n  [\"553\u00d73=\", \"466\u00d72=\"],\n  [\"281\u00d73=\", \"923\u00d72=\"],\n  [\"425\u00d78=\", \"739\u00d73=\"],\n  [\"677\u00d79=\", \"818\u00d76=\"],\n  [\"288\u00d72=\", \"382\u00d79=\"],\n  [\"413\u00d73=\", \"700\u00d73=\"],\n  [\"847\u00d77=\", \"485\u00d73=\"],\n  [\"166\u00d76=\", \"556\u00d73=\"],\n  [\"199\u00d78=\", \"343\u00d72=\"],\n  [\"974\u00d77=\", \"619\u00d72=\"],\n  [\"993\u00d76=\", \"199\u00d76=\"],\n  [\"977\u00d79=\", \"677\u00d73=\"],\n  [\"406\u00d79=\", \"187\u00d76=\"],\n  [\"577\u00d79=\", \"181\u00d78=\"],\n  [\"239\u00d74=\", \"782\u00d76=\"],\n  [\"580\u00d72=\", \"670\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication problems with new\n# operands, keeping every other part of the document (layout, fonts,\n# formatting, trailing \"=\") unchanged. Each original expression is unique\n# in the document, so a Find/Replace (ReplaceAll) per pair is safe and\n# only ever touches the single intended cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"789\u00d74=\", \"217\u00d76=\"),\n  @(\"362\u00d73=\", \"619\u00d75=\"),\n  @(\"791\u00d74=\", \"111\u00d78=\"),\n  @(\"201\u00d79=\", \"404\u00d72=\"),\n  @(\"329\u00d78=\", \"363\u00d76=\"),\n  @(\"262\u00d77=\", \"547\u00d79=\"),\n  @(\"225\u00d77=\", \"310\u00d78=\"),\n  @(\"443\u00d73=\", \"315\u00d75=\"),\n  @(\"867\u00d76=\", \"955\u00d77=\"),\n  @(\"553\u00d73=\", \"466\u00d72=\"),\n  @(\"281\u00d73=\", \"923\u00d72=\"),\n  @(\"425\u00d78=\", \"739\u00d73=\"),\n  @(\"677\u00d79=\", \"818\u00d76=\"),\n  @(\"288\u00d72=\", \"382\u00d79=\"),\n  @(\"413\u00d73=\", \"700\u00d73=\"),\n  @(\"847\u00d77=\", \"485\u00d73=\"),\n  @(\"166\u00d76=\", \"556\u00d73=\"),\n  @(\"199\u00d78=\", \"343\u00d72=\"),\n  @(\"974\u00d77=\", \"619\u00d72=\"),\n  @(\"993\u00d76=\", \"199\u00d76=\"),\n  @(\"977\u00d79=\", \"677\u00d73=\"),\n  @(\"406\u00d79=\", \"187\u00d76=\"),\n  @(\"577\u00d79=\", \"181\u00d78=\"),\n  @(\"239\u00d74=\", \"782\u00d76=\"),\n  @(\"580\u00d72=\", \"670\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
